$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update 想去人数 (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3395
$ws1.Range("F5").Value = 1578
$ws1.Range("F6").Value = 59
$ws1.Range("F7").Value = 327

# Sheet "全部类型" (all types) - update 想去人数 (want-to-go count) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3395
$ws4.Range("F5").Value = 1578
$ws4.Range("F6").Value = 59
$ws4.Range("F8").Value = 327
